# Auto-generated update of leve-profit market data cells (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1289.4166
$ws.Cells.Item(28, 9).Value = 1289.4166
$ws.Cells.Item(28, 11).Value = 1289.4166
$ws.Cells.Item(28, 13).Value = -804.4166
$ws.Cells.Item(98, 8).Value = 2246.8333
$ws.Cells.Item(98, 9).Value = 878.17645
$ws.Cells.Item(98, 10).Value = 5570.7144
$ws.Cells.Item(98, 11).Value = 878.17645
$ws.Cells.Item(98, 12).Value = 5570.7144
$ws.Cells.Item(98, 13).Value = 619.82355
$ws.Cells.Item(98, 14).Value = -8566.714400000001
$ws.Cells.Item(100, 8).Value = 1798.1428
$ws.Cells.Item(100, 9).Value = 1264.6666
$ws.Cells.Item(100, 11).Value = 1264.6666
$ws.Cells.Item(100, 13).Value = -723.6666
$ws.Cells.Item(111, 8).Value = 1146.25
$ws.Cells.Item(111, 9).Value = 1125.5
$ws.Cells.Item(111, 10).Value = 1250
$ws.Cells.Item(111, 11).Value = 3376.5
$ws.Cells.Item(111, 12).Value = 3750
$ws.Cells.Item(111, 13).Value = -309.5
$ws.Cells.Item(111, 14).Value = -9884
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(113, 8).Value = 4065.7778
$ws.Cells.Item(113, 9).Value = 2999
$ws.Cells.Item(113, 10).Value = 4599.1665
$ws.Cells.Item(113, 11).Value = 2999
$ws.Cells.Item(113, 12).Value = 4599.1665
$ws.Cells.Item(113, 13).Value = 255
$ws.Cells.Item(113, 14).Value = -11107.1665
$ws.Cells.Item(116, 8).Value = 5695.1113
$ws.Cells.Item(116, 9).Value = 4543.1665
$ws.Cells.Item(116, 11).Value = 4543.1665
$ws.Cells.Item(116, 13).Value = -1101.1665
$ws.Cells.Item(118, 8).Value = 1315.6
$ws.Cells.Item(118, 9).Value = 1315.6
$ws.Cells.Item(118, 11).Value = 3946.8
$ws.Cells.Item(118, 13).Value = -2289.8
$ws.Cells.Item(122, 8).Value = 2246.8333
$ws.Cells.Item(122, 9).Value = 878.17645
$ws.Cells.Item(122, 10).Value = 5570.7144
$ws.Cells.Item(122, 11).Value = 2634.52935
$ws.Cells.Item(122, 12).Value = 16712.1432
$ws.Cells.Item(122, 13).Value = -184.5293500000002
$ws.Cells.Item(122, 14).Value = -21612.1432
$ws.Cells.Item(132, 8).Value = 2390.5
$ws.Cells.Item(132, 9).Value = 2268.6
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 6805.799999999999
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -4275.799999999999
$ws.Cells.Item(132, 14).Value = -14060
$ws.Cells.Item(138, 8).Value = 6434.536
$ws.Cells.Item(138, 10).Value = 6527.1445
$ws.Cells.Item(138, 12).Value = 19581.4335
$ws.Cells.Item(138, 14).Value = -29861.4335
$ws.Cells.Item(112, 13).ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 928.75
$ws.Cells.Item(2, 9).Value = 521.8333
$ws.Cells.Item(2, 10).Value = 2149.5
$ws.Cells.Item(2, 11).Value = 521.8333
$ws.Cells.Item(2, 12).Value = 2149.5
$ws.Cells.Item(2, 13).Value = -408.8333
$ws.Cells.Item(2, 14).Value = -2375.5
$ws.Cells.Item(55, 8).Value = 62777.668
$ws.Cells.Item(55, 10).Value = 66874.875
$ws.Cells.Item(55, 12).Value = 66874.875
$ws.Cells.Item(55, 14).Value = -67504.875
$ws.Cells.Item(61, 8).Value = 5160.8423
$ws.Cells.Item(61, 10).Value = 7590.6
$ws.Cells.Item(61, 12).Value = 7590.6
$ws.Cells.Item(61, 14).Value = -8014.6
$ws.Cells.Item(116, 8).Value = 928.75
$ws.Cells.Item(116, 9).Value = 521.8333
$ws.Cells.Item(116, 10).Value = 2149.5
$ws.Cells.Item(116, 11).Value = 521.8333
$ws.Cells.Item(116, 12).Value = 2149.5
$ws.Cells.Item(116, 13).Value = 1772.1667
$ws.Cells.Item(116, 14).Value = -6737.5
$ws.Cells.Item(136, 8).Value = 5160.8423
$ws.Cells.Item(136, 10).Value = 7590.6
$ws.Cells.Item(136, 12).Value = 22771.8
$ws.Cells.Item(136, 14).Value = -27871.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 928.75
$ws.Cells.Item(3, 9).Value = 521.8333
$ws.Cells.Item(3, 10).Value = 2149.5
$ws.Cells.Item(3, 11).Value = 521.8333
$ws.Cells.Item(3, 12).Value = 2149.5
$ws.Cells.Item(3, 13).Value = -407.8333
$ws.Cells.Item(3, 14).Value = -2377.5
$ws.Cells.Item(82, 8).Value = 13874.75
$ws.Cells.Item(82, 9).Value = 13874.75
$ws.Cells.Item(82, 11).Value = 13874.75
$ws.Cells.Item(82, 13).Value = -13491.75
$ws.Cells.Item(85, 8).Value = 13874.75
$ws.Cells.Item(85, 9).Value = 13874.75
$ws.Cells.Item(85, 11).Value = 13874.75
$ws.Cells.Item(85, 13).Value = -12548.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 12834.417
$ws.Cells.Item(31, 9).Value = 38536.332
$ws.Cells.Item(31, 10).Value = 4267.1113
$ws.Cells.Item(31, 11).Value = 38536.332
$ws.Cells.Item(31, 12).Value = 4267.1113
$ws.Cells.Item(31, 13).Value = -38241.332
$ws.Cells.Item(31, 14).Value = -4857.1113
$ws.Cells.Item(34, 8).Value = 12834.417
$ws.Cells.Item(34, 9).Value = 38536.332
$ws.Cells.Item(34, 10).Value = 4267.1113
$ws.Cells.Item(34, 11).Value = 38536.332
$ws.Cells.Item(34, 12).Value = 4267.1113
$ws.Cells.Item(34, 13).Value = -38334.332
$ws.Cells.Item(34, 14).Value = -4671.1113
$ws.Cells.Item(99, 8).Value = 46415.125
$ws.Cells.Item(99, 9).Value = 61264.8
$ws.Cells.Item(99, 10).Value = 21665.666
$ws.Cells.Item(99, 11).Value = 61264.8
$ws.Cells.Item(99, 12).Value = 21665.666
$ws.Cells.Item(99, 13).Value = -59766.8
$ws.Cells.Item(99, 14).Value = -24661.666
$ws.Cells.Item(126, 8).Value = 46415.125
$ws.Cells.Item(126, 9).Value = 61264.8
$ws.Cells.Item(126, 10).Value = 21665.666
$ws.Cells.Item(126, 11).Value = 183794.4
$ws.Cells.Item(126, 12).Value = 64996.99800000001
$ws.Cells.Item(126, 13).Value = -181324.4
$ws.Cells.Item(126, 14).Value = -69936.99800000001
$ws.Cells.Item(134, 8).Value = 4394.8887
$ws.Cells.Item(134, 9).Value = 4365.2856
$ws.Cells.Item(134, 11).Value = 13095.8568
$ws.Cells.Item(134, 13).Value = -10560.8568
$ws.Cells.Item(141, 8).Value = 539102.7
$ws.Cells.Item(141, 10).Value = 539102.7
$ws.Cells.Item(141, 12).Value = 539102.7
$ws.Cells.Item(141, 14).Value = -549462.7

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 138686.12
$ws.Cells.Item(94, 10).Value = 19900
$ws.Cells.Item(94, 12).Value = 59700
$ws.Cells.Item(94, 14).Value = -61052
$ws.Cells.Item(112, 8).Value = 13589.375
$ws.Cells.Item(112, 9).Value = 7178.75
$ws.Cells.Item(112, 11).Value = 21536.25
$ws.Cells.Item(112, 13).Value = -20428.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 1585.5714
$ws.Cells.Item(22, 9).Value = 1279.8
$ws.Cells.Item(22, 11).Value = 1279.8
$ws.Cells.Item(22, 13).Value = -750.8
$ws.Cells.Item(96, 8).Value = 19985.5
$ws.Cells.Item(96, 10).Value = 19985.5
$ws.Cells.Item(96, 12).Value = 19985.5
$ws.Cells.Item(96, 14).Value = -25477.5
$ws.Cells.Item(97, 8).Value = 703
$ws.Cells.Item(97, 9).Value = 621.875
$ws.Cells.Item(97, 10).Value = 1027.5
$ws.Cells.Item(97, 11).Value = 621.875
$ws.Cells.Item(97, 12).Value = 1027.5
$ws.Cells.Item(97, 13).Value = -125.875
$ws.Cells.Item(97, 14).Value = -2019.5
$ws.Cells.Item(107, 8).Value = 519.0909
$ws.Cells.Item(107, 9).Value = 473
$ws.Cells.Item(107, 11).Value = 473
$ws.Cells.Item(107, 13).Value = 1447

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 111111
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(93, 8).Value = 866.75
$ws.Cells.Item(93, 9).Value = 695
$ws.Cells.Item(93, 11).Value = 695
$ws.Cells.Item(93, 13).Value = 553
$ws.Cells.Item(113, 8).Value = 111111
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(122, 8).Value = 6676.7104
$ws.Cells.Item(122, 9).Value = 6307.4443
$ws.Cells.Item(122, 10).Value = 7583.091
$ws.Cells.Item(122, 11).Value = 18922.3329
$ws.Cells.Item(122, 12).Value = 22749.273
$ws.Cells.Item(122, 13).Value = -16472.3329
$ws.Cells.Item(122, 14).Value = -27649.273
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(113, 14).ClearContents()
